# Fill in the remaining consolidated metrics for the first three employees
# on the "Consolidated Data" sheet, and make that sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consolidated Data")

# Row 2 - Chrissy Cummings
# "4.6" looks numeric, so force text storage (and then restore the default
# style so no extra formatting is left behind on the cell).
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.6"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 4166
$ws.Range("G2").Value = 4488.02
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 29.5
$ws.Range("K2").Value = 31.78
$ws.Range("L2").Value = 4195.5
$ws.Range("M2").Value = 4519.8

# Row 3 - Danielle Mai
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.9"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 6266.85
$ws.Range("G3").Value = 6750.95
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 21.54
$ws.Range("L3").Value = 6286.85
$ws.Range("M3").Value = 6772.49

# Row 4 - Jasmine Saiz
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = 596
$ws.Range("G4").Value = 641.8
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 32.31
$ws.Range("L4").Value = 626
$ws.Range("M4").Value = 674.1099999999999

# Make the "Consolidated Data" sheet the active tab
$ws.Activate()
